# Update "想去人数" (F column) values on the 展览, 演出, and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3541
$ws1.Range("F4").Value = 150
$ws1.Range("F5").Value = 7035
$ws1.Range("F6").Value = 3504
$ws1.Range("F7").Value = 61
$ws1.Range("F8").Value = 157
$ws1.Range("F11").Value = 90
$ws1.Range("F12").Value = 50
$ws1.Range("F13").Value = 33
$ws1.Range("F14").Value = 184
$ws1.Range("F15").Value = 599
$ws1.Range("F16").Value = 50
$ws1.Range("F17").Value = 47

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 29

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3541
$ws4.Range("F3").Value = 29
$ws4.Range("F5").Value = 150
$ws4.Range("F6").Value = 7035
$ws4.Range("F7").Value = 3504
$ws4.Range("F8").Value = 61
$ws4.Range("F9").Value = 157
$ws4.Range("F12").Value = 90
$ws4.Range("F13").Value = 50
$ws4.Range("F14").Value = 33
$ws4.Range("F15").Value = 184
$ws4.Range("F16").Value = 599
$ws4.Range("F17").Value = 50
$ws4.Range("F18").Value = 47
